$wb = $excel.ActiveWorkbook

# --- profile sheet: insert new column for USDA soil order ---
$wsProfile = $wb.Worksheets.Item("profile")
$wsProfile.Columns.Item(14).Insert()
$wsProfile.Cells.Item(1, 14).Value = "pro_usda_soil_order"
$wsProfile.Cells.Item(4, 14).Value = "Gelisols"

# restore profile sheet's view state
$wsProfile.Activate()
$wsProfile.Range("O11").Select()

# --- metadata sheet becomes the active tab/selection ---
$wsMeta = $wb.Worksheets.Item("metadata")
$wsMeta.Activate()
$wsMeta.Range("A4").Select()
